$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("E2").Value = '2023.12.22 10:00 - 2024.05.05 20:00'
$ws.Range("E3").Value = '2024.03.08 10:00 - 04.06 20:30'
$ws.Range("E4").Value = '2024.03.16 10:00 - 03.17 17:00'
$ws.Range("F4").Value = 7894
$ws.Range("E5").Value = '2024.03.16 10:00 - 03.17 17:00'
$ws.Range("F5").Value = 7894
$ws.Range("E6").Value = '2024.03.16 10:00 - 03.16 17:00'
$ws.Range("E7").Value = '2024.03.16 10:00 - 03.16 17:00'
$ws.Range("E8").Value = '2024.03.16 10:00 - 03.17 17:00'
$ws.Range("F8").Value = 2137
$ws.Range("E9").Value = '2024.03.16 10:00 - 03.17 17:00'
$ws.Range("F9").Value = 8650
$ws.Range("E10").Value = '2024.03.16 11:00 - 03.16 19:00'
$ws.Range("E11").Value = '2024.03.17 10:00 - 03.17 17:00'
$ws.Range("E12").Value = '2024.03.17 10:00 - 03.17 17:00'
$ws.Range("E13").Value = '2024.03.23 10:00 - 03.24 17:00'
$ws.Range("F13").Value = 5788
$ws.Range("E14").Value = '2024.03.23 10:00 - 03.23 17:00'
$ws.Range("E15").Value = '2024.03.23 10:00 - 03.24 17:00'
$ws.Range("F15").Value = 2782
$ws.Range("E16").Value = '2024.03.23 10:00 - 03.23 17:00'
$ws.Range("E17").Value = '2024.03.24 10:30 - 03.24 18:00'
$ws.Range("F17").Value = 419
$ws.Range("E18").Value = '2024.03.24 10:00 - 03.24 17:00'
$ws.Range("E19").Value = '2024.03.29 14:00 - 03.31 20:00'
$ws.Range("E20").Value = '2024.03.30 10:00 - 03.30 17:00'
$ws.Range("E21").Value = '2024.03.31 10:00 - 03.31 16:00'
$ws.Range("F21").Value = 97
$ws.Range("E22").Value = '2024.04.04 09:30 - 04.05 17:00'
$ws.Range("F22").Value = 3961
$ws.Range("E23").Value = '2024.04.04 10:30 - 04.04 13:30'
$ws.Range("E24").Value = '2024.04.04 09:30 - 04.05 17:00'
$ws.Range("F24").Value = 65
$ws.Range("E25").Value = '2024.04.04 09:00 - 04.05 17:00'
$ws.Range("E26").Value = '2024.04.04 09:30 - 04.05 17:00'
$ws.Range("E27").Value = '2024.04.04 11:35 - 04.04 14:50'
$ws.Range("E28").Value = '2024.04.04 09:30 - 04.05 17:00'
$ws.Range("F28").Value = 30
$ws.Range("E29").Value = '2024.04.04 09:00 - 04.05 17:00'
$ws.Range("F29").Value = 5629
$ws.Range("E30").Value = '2024.04.04 10:00 - 04.06 17:00'
$ws.Range("F30").Value = 213
$ws.Range("E31").Value = '2024.04.05 10:30 - 04.05 13:45'
$ws.Range("F31").Value = 75
$ws.Range("E32").Value = '2024.04.06 10:00 - 04.06 17:00'
$ws.Range("E33").Value = '2024.04.13 10:00 - 04.14 17:00'
$ws.Range("E34").Value = '2024.04.13 10:00 - 04.13 17:00'
$ws.Range("F34").Value = 406
$ws.Range("E35").Value = '2024.04.19 10:00 - 04.21 17:00'
$ws.Range("F35").Value = 2954
$ws.Range("E36").Value = '2024.04.20 09:00 - 04.21 17:00'
$ws.Range("F36").Value = 1537
$ws.Range("E37").Value = '2024.04.20 10:30 - 04.20 17:00'
$ws.Range("E38").Value = '2024.04.20 09:00 - 04.21 17:00'
$ws.Range("F38").Value = 1239
$ws.Range("E39").Value = '2024.05.01 09:30 - 05.03 17:00'
$ws.Range("F39").Value = 5159
$ws.Range("E40").Value = '2024.05.01 09:30 - 05.03 17:00'
$ws.Range("E41").Value = '2024.05.01 09:30 - 05.03 17:00'
$ws.Range("E42").Value = '2024.05.01 09:00 - 05.04 17:00'
$ws.Range("E43").Value = '2024.05.01 09:00 - 05.04 17:00'
$ws.Range("F43").Value = 3640
$ws.Range("E44").Value = '2024.05.03 10:30 - 05.03 15:00'
$ws.Range("E45").Value = '2024.05.04 11:00 - 05.04 18:00'
$ws.Range("E46").Value = '2024.05.18 10:00 - 05.19 17:00'
$ws.Range("E47").Value = '2024.05.18 10:00 - 05.18 17:00'
$ws.Range("E48").Value = '2024.05.18 10:00 - 05.18 17:00'
$ws.Range("E49").Value = '2024.05.25 09:00 - 05.26 17:00'
$ws.Range("E50").Value = '2024.05.25 10:00 - 05.25 17:00'

$ws = $wb.Worksheets.Item("演出")
$ws.Range("E2").Value = '2024.03.27 19:30 - 03.27 21:10'
$ws.Range("E3").Value = '2024.03.30 14:30 - 03.30 16:00'
$ws.Range("F3").Value = 157
$ws.Range("E4").Value = '2024.03.30 18:00 - 03.30 22:00'
$ws.Range("F4").Value = 18
$ws.Range("E5").Value = '2024.04.14 19:30 - 04.14 21:00'
$ws.Range("E6").Value = '2024.04.20 19:30 - 04.20 21:00'
$ws.Range("E7").Value = '2024.04.27 20:00 - 04.27 21:30'
$ws.Range("E8").Value = '2024.05.11 19:30 - 05.11 21:00'
$ws.Range("E9").Value = '2024.05.18 20:00 - 05.18 22:00'
$ws.Range("E10").Value = '2024.05.25 19:30 - 06.06 22:00'
$ws.Range("E11").Value = '2024.06.28 19:30 - 06.28 21:00'

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("E2").Value = '2023.12.22 10:00 - 2024.03.15 17:00'
$ws.Range("E3").Value = '2023.12.29 00:00 - 2024.03.31 23:59'
$ws.Range("F3").Value = 1372

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("E2").Value = '2023.12.22 10:00 - 2024.03.15 17:00'
$ws.Range("E3").Value = '2023.12.29 00:00 - 2024.03.31 23:59'
$ws.Range("F3").Value = 1372
$ws.Range("E4").Value = '2024.03.08 10:00 - 04.06 20:30'
$ws.Range("E5").Value = '2024.03.16 10:00 - 03.17 17:00'
$ws.Range("F5").Value = 7894
$ws.Range("E6").Value = '2024.03.16 10:00 - 03.17 17:00'
$ws.Range("F6").Value = 7894
$ws.Range("E7").Value = '2024.03.16 10:00 - 03.16 17:00'
$ws.Range("E8").Value = '2024.03.16 10:00 - 03.16 17:00'
$ws.Range("E9").Value = '2024.03.16 10:00 - 03.17 17:00'
$ws.Range("F9").Value = 2137
$ws.Range("E10").Value = '2024.03.16 10:00 - 03.17 17:00'
$ws.Range("F10").Value = 8650
$ws.Range("E11").Value = '2024.03.16 11:00 - 03.16 19:00'
$ws.Range("E12").Value = '2024.03.17 10:00 - 03.17 17:00'
$ws.Range("E13").Value = '2024.03.17 10:00 - 03.17 17:00'
$ws.Range("E14").Value = '2024.03.23 10:00 - 03.24 17:00'
$ws.Range("F14").Value = 5788
$ws.Range("E15").Value = '2024.03.23 10:00 - 03.23 17:00'
$ws.Range("E16").Value = '2024.03.23 10:00 - 03.24 17:00'
$ws.Range("F16").Value = 2782
$ws.Range("E17").Value = '2024.03.23 10:00 - 03.23 17:00'
$ws.Range("E18").Value = '2024.03.24 10:30 - 03.24 18:00'
$ws.Range("F18").Value = 419
$ws.Range("E19").Value = '2024.03.24 10:00 - 03.24 17:00'
$ws.Range("E20").Value = '2024.03.29 14:00 - 03.31 20:00'
$ws.Range("E21").Value = '2024.03.30 14:30 - 03.30 16:00'
$ws.Range("F21").Value = 157
$ws.Range("E22").Value = '2024.03.30 10:00 - 03.30 17:00'
$ws.Range("E23").Value = '2024.03.31 10:00 - 03.31 16:00'
$ws.Range("F23").Value = 97
$ws.Range("E24").Value = '2024.04.04 09:30 - 04.05 17:00'
$ws.Range("F24").Value = 3961
$ws.Range("E25").Value = '2024.04.04 10:30 - 04.04 13:30'
$ws.Range("E26").Value = '2024.04.04 09:30 - 04.05 17:00'
$ws.Range("F26").Value = 65
$ws.Range("E27").Value = '2024.04.04 09:30 - 04.05 17:00'
$ws.Range("E28").Value = '2024.04.04 11:35 - 04.04 14:50'
$ws.Range("E29").Value = '2024.04.04 09:30 - 04.05 17:00'
$ws.Range("F29").Value = 30
$ws.Range("E30").Value = '2024.04.04 09:00 - 04.05 17:00'
$ws.Range("F30").Value = 5629
$ws.Range("E31").Value = '2024.04.05 10:30 - 04.05 13:45'
$ws.Range("F31").Value = 75
$ws.Range("E32").Value = '2024.04.06 10:00 - 04.06 17:00'
$ws.Range("E33").Value = '2024.04.13 10:00 - 04.14 17:00'
$ws.Range("E34").Value = '2024.04.13 10:00 - 04.13 17:00'
$ws.Range("F34").Value = 406
$ws.Range("E35").Value = '2024.04.19 10:00 - 04.21 17:00'
$ws.Range("F35").Value = 2954
$ws.Range("E36").Value = '2024.04.20 09:00 - 04.21 17:00'
$ws.Range("F36").Value = 1537
$ws.Range("E37").Value = '2024.04.20 19:30 - 04.20 21:00'
$ws.Range("E38").Value = '2024.04.20 10:30 - 04.20 17:00'
$ws.Range("E39").Value = '2024.04.20 09:00 - 04.21 17:00'
$ws.Range("F39").Value = 1249
$ws.Range("E40").Value = '2024.04.27 20:00 - 04.27 21:30'
$ws.Range("E41").Value = '2024.05.01 09:30 - 05.03 17:00'
$ws.Range("F41").Value = 5159
$ws.Range("E42").Value = '2024.05.01 09:30 - 05.03 17:00'
$ws.Range("E43").Value = '2024.05.01 09:30 - 05.03 17:00'
$ws.Range("E44").Value = '2024.05.01 09:00 - 05.04 17:00'
$ws.Range("F44").Value = 3640
$ws.Range("E45").Value = '2024.05.04 11:00 - 05.04 18:00'
$ws.Range("E46").Value = '2024.05.18 20:00 - 05.18 22:00'
$ws.Range("E47").Value = '2024.05.18 10:00 - 05.19 17:00'
$ws.Range("E48").Value = '2024.05.18 10:00 - 05.18 17:00'
$ws.Range("E49").Value = '2024.05.25 09:00 - 05.26 17:00'
$ws.Range("E50").Value = '2024.06.28 19:30 - 06.28 21:00'
